## Update the cci-ecv vocabulary spreadsheet with changes from the
## python3_updates branch: two new climate-data-record rows
## (vegetation parameters, river discharge) and a couple of small
## workbook bookkeeping tweaks that came along with the resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the two new vocabulary rows (row 23 + row 24) ----------------

$ws.Range("A23").Value = "cciecv_vegParam"
$ws.Range("B23").Value = "vegetation parameters"
$ws.Range("C23").Value = "VEGETATION"
$ws.Range("D23").Value = "Vegetation Parameters climate data record produced from satellite data as part of the European Space Agency (ESA) Climate Change Initiative (CCI)"

$ws.Range("A24").Value = "cciecv_riverDischarge"
$ws.Range("B24").Value = "river discharge"
$ws.Range("C24").Value = "RD"
$ws.Range("D24").Value = "River discharge climate data record produced from satellite data as part of the European Space Agency (ESA) Climate Change Initiative (CCI)"

# --- Drop the unused "Heading1" cell style ---------------------------------
# (left over from an earlier template; pruned during this resave along with
# its now-orphaned font)

$wb.Styles("Heading1").Delete()

# --- Turn on iterative calculation with a tight convergence delta ---------

$excel.Iteration = $true
$excel.MaxChange = 0.0001
